# Updated cryptos list values (Price and Volume(1h) columns) refreshed by GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.220.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.380.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.379.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.119.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.376.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.95%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '558.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.497.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.48%  '
$ws.Range("E32").Value = '  -2.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("E36").Value = '  +5.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.72%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("E45").Value = '  +4.90%  '
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
